$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @("D2", "E2", "D3", "E3", "E4", "D5", "E5", "D6", "E6", "E7", "E8", "E9", "E10", "E11", "D12", "E13", "D14", "E14", "D15", "E15", "E16", "D17", "E17", "D18", "E18", "D19", "E19", "E20", "D21", "E21", "E22", "D23", "E24", "D25", "E25", "E26", "E27", "E28", "E29", "D30", "E30", "E31", "D32", "E32", "E33", "D34", "E34", "E35", "E36", "E37", "D38", "E38", "D39", "E39", "E40", "D41", "D42", "E42", "E43", "D44", "E44", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "E50", "E51")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '62.028.42'
$ws.Range("E2").Value = '  -2.00%  '
$ws.Range("D3").Value = '3.418.66'
$ws.Range("E3").Value = '  -1.39%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '578.12'
$ws.Range("E5").Value = '  -0.59%  '
$ws.Range("D6").Value = '153.14'
$ws.Range("E6").Value = '  +3.78%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +1.28%  '
$ws.Range("E9").Value = '  +3.75%  '
$ws.Range("E10").Value = '  -0.56%  '
$ws.Range("E11").Value = '  +3.09%  '
$ws.Range("D12").Value = '4.002.82'
$ws.Range("E13").Value = '  +0.62%  '
$ws.Range("D14").Value = '28.69'
$ws.Range("E14").Value = '  -2.76%  '
$ws.Range("D15").Value = '3.421.12'
$ws.Range("E15").Value = '  -1.54%  '
$ws.Range("E16").Value = '  -0.39%  '
$ws.Range("D17").Value = '62.041.48'
$ws.Range("E17").Value = '  -2.02%  '
$ws.Range("D18").Value = '6.51'
$ws.Range("E18").Value = '  +1.79%  '
$ws.Range("D19").Value = '14.51'
$ws.Range("E19").Value = '  +0.02%  '
$ws.Range("E20").Value = '  -4.11%  '
$ws.Range("D21").Value = '381.97'
$ws.Range("E21").Value = '  -1.72%  '
$ws.Range("E22").Value = '  +1.09%  '
$ws.Range("D23").Value = '75.24'
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").Value = '3.560.89'
$ws.Range("E25").Value = '  -1.50%  '
$ws.Range("E26").Value = '  -3.57%  '
$ws.Range("E27").Value = '  -1.54%  '
$ws.Range("E28").Value = '  +0.15%  '
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("D30").Value = '7.91'
$ws.Range("E30").Value = '  -3.71%  '
$ws.Range("E31").Value = '  -1.07%  '
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").Value = '  -0.04%  '
$ws.Range("E33").Value = '  -0.76%  '
$ws.Range("D34").Value = '23.23'
$ws.Range("E34").Value = '  -1.01%  '
$ws.Range("E35").Value = '  +3.60%  '
$ws.Range("E36").Value = '  -1.18%  '
$ws.Range("E37").Value = '  -2.56%  '
$ws.Range("D38").Value = '168.44'
$ws.Range("E38").Value = '  +0.16%  '
$ws.Range("D39").Value = '30.92'
$ws.Range("E39").Value = '  -3.45%  '
$ws.Range("E40").Value = '  -1.57%  '
$ws.Range("D41").Value = '0.0786'
$ws.Range("D42").Value = '42.69'
$ws.Range("E42").Value = '  +0.66%  '
$ws.Range("E43").Value = '  -1.73%  '
$ws.Range("D44").Value = '4.41'
$ws.Range("E44").Value = '  +0.70%  '
$ws.Range("E45").Value = '  -3.88%  '
$ws.Range("D46").Value = '1.16'
$ws.Range("E46").Value = '  -3.57%  '
$ws.Range("D47").Value = '2.552.18'
$ws.Range("E47").Value = '  -1.63%  '
$ws.Range("D48").Value = '6.86'
$ws.Range("E48").Value = '  +0.46%  '
$ws.Range("D49").Value = '22.62'
$ws.Range("E49").Value = '  -1.70%  '
$ws.Range("E50").Value = '  -5.65%  '
$ws.Range("E51").Value = '  -0.04%  '
